$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# Row 1: add headers for Trial 2 / Trial 3 / Trial 4 columns (D1,E1,F1)
# Copy format (fill/border/font) from C1, which already carries the
# "Trial" header style, then set the text.
# ---------------------------------------------------------------------
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1:F1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("D1").Value = "Trial 2"
$ws.Range("E1").Value = "Trial 3"
$ws.Range("F1").Value = "Trial 4"

# ---------------------------------------------------------------------
# Row 2: Date row - update trial 1 date, fill trial 2/3/4 dates
# ---------------------------------------------------------------------
$ws.Range("C2").Value = 42131
$ws.Range("D2").Value = 42131
$ws.Range("E2").Value = 42131
$ws.Range("F2").Value = 42131

# ---------------------------------------------------------------------
# Row 3: Hour row (the literal text for E3 is filled in further below,
# to match the order in which new shared strings were introduced)
# ---------------------------------------------------------------------
$ws.Range("C3").Value = 0.49583333333333335
$ws.Range("D3").NumberFormat = "h:mm AM/PM"
$ws.Range("D3").Value = 0.4993055555555555
$ws.Range("F3").NumberFormat = "h:mm AM/PM"
$ws.Range("F3").Value = 0.69652777777777775

# ---------------------------------------------------------------------
# Row 5: ntrials - fill the whole row (C through M) with 500
# ---------------------------------------------------------------------
$ws.Range("C5:M5").NumberFormat = "General"
$ws.Range("C5:M5").Value = 500

# ---------------------------------------------------------------------
# Row 6: cBF.N
# ---------------------------------------------------------------------
$ws.Range("C6").Value = 6
$ws.Range("D6:E6").NumberFormat = "General"
$ws.Range("D6:E6").Value = 6

# ---------------------------------------------------------------------
# Row 7: aBF.N
# ---------------------------------------------------------------------
$ws.Range("C7").Value = 6
$ws.Range("D7:E7").NumberFormat = "General"
$ws.Range("D7:E7").Value = 6

# ---------------------------------------------------------------------
# Row 8: cBF.a
# ---------------------------------------------------------------------
$ws.Range("C8:E8").NumberFormat = "0.000000"
$ws.Range("C8:E8").Value = 0.00005

# ---------------------------------------------------------------------
# Row 9: aBF.a
# ---------------------------------------------------------------------
$ws.Range("C9:E9").NumberFormat = "0.000000"
$ws.Range("C9:E9").Value = 0.00005

# ---------------------------------------------------------------------
# Row 10: cBF.rn
# ---------------------------------------------------------------------
$ws.Range("D10:E10").Value = "[-1 1;-1 1]"

# ---------------------------------------------------------------------
# Row 11: aBF.rn
# ---------------------------------------------------------------------
$ws.Range("D11:E11").Value = "[-1 1;-1 1]"

# ---------------------------------------------------------------------
# Row 12: cBF.T
# ---------------------------------------------------------------------
$ws.Range("C12:E12").NumberFormat = "0.000"
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 1.3
$ws.Range("E12").Value = 2

# ---------------------------------------------------------------------
# Row 13: aBF.T
# ---------------------------------------------------------------------
$ws.Range("C13:E13").NumberFormat = "0.000"
$ws.Range("C13").Value = 3
$ws.Range("D13").Value = 3
$ws.Range("E13").Value = 2

# ---------------------------------------------------------------------
# Row 14: par.Q
# ---------------------------------------------------------------------
$ws.Range("C14:E14").NumberFormat = "0.00E+00"
$ws.Range("C14:E14").Value = 1000

# ---------------------------------------------------------------------
# Row 15: par.R
# ---------------------------------------------------------------------
$ws.Range("C15:E15").NumberFormat = "General"
$ws.Range("C15:E15").Value = 100

# ---------------------------------------------------------------------
# Row 16: par.uSat
# ---------------------------------------------------------------------
$ws.Range("C16:E16").NumberFormat = "0.00000"
$ws.Range("C16:E16").Value = 0.009

# ---------------------------------------------------------------------
# Row 17: par.varRand
# ---------------------------------------------------------------------
$ws.Range("C17:E17").NumberFormat = "0.00000"
$ws.Range("C17:E17").Value = 0.0001

# ---------------------------------------------------------------------
# Row 18: par.expSteps
# ---------------------------------------------------------------------
$ws.Range("C18:E18").NumberFormat = "0.00000"
$ws.Range("C18:E18").Value = 1

# ---------------------------------------------------------------------
# Row 19: par.gamma
# ---------------------------------------------------------------------
$ws.Range("C19:E19").NumberFormat = "0.00000"
$ws.Range("C19:E19").Value = 0.97

# ---------------------------------------------------------------------
# Row 20: par.cost
# ---------------------------------------------------------------------
$ws.Range("D20:F20").Value = "quadratic"

# ---------------------------------------------------------------------
# Row 23: converge?   (C23 introduces the "no" string)
# ---------------------------------------------------------------------
$ws.Range("C23").Value = "no"
$ws.Range("D23").Value = "no"

# ---------------------------------------------------------------------
# Row 24: performance (C24 introduces "bad", D24 introduces "deviating")
# ---------------------------------------------------------------------
$ws.Range("C24").Value = "bad"
$ws.Range("D24").Value = "deviating"

# ---------------------------------------------------------------------
# Row 3 (cont'd): E3 introduces the "14:42 PM" literal text string
# ---------------------------------------------------------------------
$ws.Range("E3").Value = "14:42 PM"

# ---------------------------------------------------------------------
# Row 26: remark (E26 introduces the "stopped..." string)
# ---------------------------------------------------------------------
$ws.Range("E26").Value = "stopped at around 250th trial due to instability induced by sudden leap in the alpha. This can be used for further tuning"

# ---------------------------------------------------------------------
# Row 23 (cont'd): E23 introduces the "yes" string
# ---------------------------------------------------------------------
$ws.Range("E23").Value = "yes"

# ---------------------------------------------------------------------
# Row 24 (cont'd): E24 introduces the "so so" string
# ---------------------------------------------------------------------
$ws.Range("E24").Value = "so so"

# ---------------------------------------------------------------------
# Row 25: cost fun
# ---------------------------------------------------------------------
$ws.Range("C25").Value = "quadratic"
$ws.Range("D25").Value = "quadratic"
$ws.Range("E25").Value = "quadratic"

# ---------------------------------------------------------------------
# Row 27: label changes from "conv steps" to "alfa filter coeff"
# (B27 introduces "alfa filter coeff", E27 introduces "1/300?? (fail!)",
# F27 introduces "x10")
# ---------------------------------------------------------------------
$ws.Range("B27").Value = "alfa filter coeff"
$ws.Range("E27").Value = "1/300?? (fail!)"
$ws.Range("F27").Value = "x10"

# ---------------------------------------------------------------------
# Update selection to match the author's final cursor position
# ---------------------------------------------------------------------
$ws.Range("E19").Select() | Out-Null
